$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" property value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-04-12T14:29:05+00:00"

# --- 2. "Mapping Table 3" (Quantité_composant_prescrite mapping):
#        collapse the detailed Nombre/Unité sub-rows into a single,
#        not-related-to summary row ---
$ws = $wb.Worksheets.Item("Mapping Table 3")

# Update row 3: rename the source path (drop the "/Nombre" suffix),
# mark the relationship as not-related-to, and clear the target column.
$ws.Cells.Item(3, 1).Value = "Messages/M_prescription_médicaments/Prescription/Elément_prescr_médic/Composant_prescrit/Quantité_composant_prescrite"
$ws.Cells.Item(3, 3).Value = "not-related-to"
$ws.Cells.Item(3, 4).Value = ""

# Remove the now-redundant detail rows (old rows 4-6: Unité/value/code/unit).
$ws.Rows("4:6").Delete()
